$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Insert a new bold "DATA ANALYSIS" heading paragraph before the
#    existing opening paragraph, and trim the trailing clause from that
#    opening paragraph.
# ---------------------------------------------------------------------
$firstPara = $d.Paragraphs(1)
$firstPara.Range.InsertParagraphBefore()

$headingPara = $d.Paragraphs(1)
$headingRange = $headingPara.Range
$headingRange.Text = "DATA ANALYSIS"
$headingRange.Font.Bold = $true
$headingRange.Font.BoldBi = $true
$headingPara.LeftIndent = 72
$headingPara.FirstLineIndent = 36

$d.Content.Find.Execute("since the data was pretty interesting", $false, $false, $false, $false, $false, $true, 1, $false, "", 1)

# ---------------------------------------------------------------------
# 2. Fix the "os" -> "is" typo and append a new closing sentence to the
#    "Score of budget range..." paragraph.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("score os inversely", $false, $false, $false, $false, $false, $true, 1, $false, "score is inversely", 2)
$d.Content.Find.Execute("but this is not the case", $false, $false, $false, $false, $false, $true, 1, $false, "but this is not the case. Also note that overall, the reading scores are higher than math scores irrespective of the per student budget. ", 2)

# ---------------------------------------------------------------------
# 3. Append a sentence about reading/math score gaps to the paragraph
#    discussing the school-type plot.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Next we look at the budget per school type", $false, $false, $false, $false, $false, $true, 1, $false, "Next we look at the budget per school type. Note that the reading and math scores are very close for charter schools but for district schools, the gap between reading and math scores widens.", 2)

# ---------------------------------------------------------------------
# 4. Rework the back half of the CONCLUSION paragraph.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("have special needs and not based on merit. So need more budget and still cant get higher scores. ", $false, $false, $false, $false, $false, $true, 1, $false, "also cover special  need students and hence get more budget allocation. They are not based on merit. ", 2)
